$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 722.3333
$ws.Range("I19").Value = 546.7
$ws.Range("J19").Value = 882
$ws.Range("K19").Value = 546.7
$ws.Range("L19").Value = 882
$ws.Range("M19").Value = -371.7
$ws.Range("N19").Value = -1232

$ws.Range("H41").Value = 606.2778
$ws.Range("I41").Value = 390.0909
$ws.Range("J41").Value = 946
$ws.Range("K41").Value = 390.0909
$ws.Range("L41").Value = 946
$ws.Range("M41").Value = 49.90910000000002
$ws.Range("N41").Value = -1826

$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()

$ws.Range("H98").Value = 1618.0605
$ws.Range("I98").Value = 1237
$ws.Range("J98").Value = 2494.5
$ws.Range("K98").Value = 1237
$ws.Range("L98").Value = 2494.5
$ws.Range("M98").Value = 261
$ws.Range("N98").Value = -5490.5

$ws.Range("H113").Value = 3896.8235
$ws.Range("I113").Value = 3182.5
$ws.Range("J113").Value = 4531.778
$ws.Range("K113").Value = 3182.5
$ws.Range("L113").Value = 4531.778
$ws.Range("M113").Value = 71.5
$ws.Range("N113").Value = -11039.778

$ws.Range("H116").Value = 3178.1
$ws.Range("I116").Value = 2352.7058
$ws.Range("J116").Value = 4257.4614
$ws.Range("K116").Value = 2352.7058
$ws.Range("L116").Value = 4257.4614
$ws.Range("M116").Value = 1089.2942
$ws.Range("N116").Value = -11141.4614

$ws.Range("H122").Value = 1618.0605
$ws.Range("I122").Value = 1237
$ws.Range("J122").Value = 2494.5
$ws.Range("K122").Value = 3711
$ws.Range("L122").Value = 7483.5
$ws.Range("M122").Value = -1261
$ws.Range("N122").Value = -12383.5

$ws.Range("H137").Value = 2461.9333
$ws.Range("I137").Value = 2769.889
$ws.Range("J137").Value = 2000
$ws.Range("K137").Value = 8309.667000000001
$ws.Range("L137").Value = 6000
$ws.Range("M137").Value = -5759.667000000001

$ws.Range("H138").Value = 1687.1666
$ws.Range("I138").Value = 1327.7037
$ws.Range("J138").Value = 2765.5557
$ws.Range("K138").Value = 3983.1111
$ws.Range("L138").Value = 8296.667099999999
$ws.Range("M138").Value = 1156.8889
$ws.Range("N138").Value = -18576.6671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6226.2554
$ws.Range("I32").Value = 6698.5713
$ws.Range("J32").Value = 2258.8
$ws.Range("K32").Value = 6698.5713
$ws.Range("L32").Value = 2258.8
$ws.Range("M32").Value = -6411.5713

$ws.Range("H132").Value = 7511
$ws.Range("I132").Value = 3883.2
$ws.Range("J132").Value = 8578
$ws.Range("K132").Value = 11649.6
$ws.Range("L132").Value = 25734
$ws.Range("M132").Value = -9119.599999999999
$ws.Range("N132").Value = -30794

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2433.1836
$ws.Range("I31").Value = 1775.037
$ws.Range("J31").Value = 3240.9092
$ws.Range("K31").Value = 1775.037
$ws.Range("L31").Value = 3240.9092
$ws.Range("M31").Value = -1480.037
$ws.Range("N31").Value = -3830.9092

$ws.Range("H34").Value = 2433.1836
$ws.Range("I34").Value = 1775.037
$ws.Range("J34").Value = 3240.9092
$ws.Range("K34").Value = 1775.037
$ws.Range("L34").Value = 3240.9092
$ws.Range("M34").Value = -1573.037
$ws.Range("N34").Value = -3644.9092

$ws.Range("H58").Value = 1566972.8
$ws.Range("I58").Value = 3756.4614
$ws.Range("J58").Value = 2636542
$ws.Range("K58").Value = 3756.4614
$ws.Range("L58").Value = 2636542
$ws.Range("M58").Value = -3553.4614
$ws.Range("N58").Value = -2636948

$ws.Range("H107").Value = 1559.3125
$ws.Range("I107").Value = 977.1818
$ws.Range("J107").Value = 2840
$ws.Range("K107").Value = 977.1818
$ws.Range("L107").Value = 2840
$ws.Range("M107").Value = 942.8182
$ws.Range("N107").Value = -6680

$ws.Range("H132").Value = 3031
$ws.Range("I132").Value = 2034
$ws.Range("J132").Value = 3404.875
$ws.Range("K132").Value = 6102
$ws.Range("L132").Value = 10214.625
$ws.Range("M132").Value = -3572
$ws.Range("N132").Value = -15274.625

$ws.Range("H134").Value = 3068.2632
$ws.Range("I134").Value = 1573.75
$ws.Range("J134").Value = 4155.1816
$ws.Range("K134").Value = 4721.25
$ws.Range("L134").Value = 12465.5448
$ws.Range("M134").Value = -2186.25
$ws.Range("N134").Value = -17535.5448

$ws.Range("H136").Value = 1566972.8
$ws.Range("I136").Value = 3756.4614
$ws.Range("J136").Value = 2636542
$ws.Range("K136").Value = 11269.3842
$ws.Range("L136").Value = 7909626
$ws.Range("M136").Value = -8719.3842
$ws.Range("N136").Value = -7914726

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 898.65515
$ws.Range("I131").Value = 525.1818
$ws.Range("J131").Value = 1126.8889
$ws.Range("K131").Value = 1575.5454
$ws.Range("L131").Value = 3380.6667
$ws.Range("M131").Value = 3464.4546
$ws.Range("N131").Value = -13460.6667

$ws.Range("H133").Value = 2600
$ws.Range("I133").Value = 2600
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 7800
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -2740
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2585.5833
$ws.Range("I16").Value = 2402.7
$ws.Range("J16").Value = 3500
$ws.Range("K16").Value = 2402.7
$ws.Range("L16").Value = 3500
$ws.Range("M16").Value = -2232.7
$ws.Range("N16").Value = -3840

$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H58").Value = 4666.6665
$ws.Range("I58").Value = 5000
$ws.Range("J58").Value = 4500
$ws.Range("K58").Value = 5000
$ws.Range("L58").Value = 4500
$ws.Range("M58").Value = -4740
$ws.Range("N58").Value = -5020

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 88215
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 88215
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 88215
$ws.Range("N46").Value = -88677

$ws.Range("H76").Value = 15000
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 15000
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 15000
$ws.Range("N76").Value = -15630

$ws.Range("H79").Value = 15000
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 15000
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 15000
$ws.Range("N79").Value = -17184

$ws.Range("H134").Value = 88215
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 88215
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 264645
$ws.Range("N134").Value = -269715

$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
